# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 1 de Mayo de 2020 a las 12:22"

# --- Rows whose country label moved (rank changed because of new totals) ---

# Row 32 used to be Japon, now shows Bielorrusia with its refreshed numbers
$ws.Range("A32").Value = "Bielorrusia"
$ws.Range("B32").Value = 14917
$ws.Range("C32").Value = 890
$ws.Range("D32").Value = 2918
$ws.Range("E32").Value = 11906
$ws.Range("F32").Value = 92
$ws.Range("G32").Value = 4
$ws.Range("H32").Value = 93

# Row 33 used to be Bielorrusia, now shows Japon (its previously unchanged numbers)
$ws.Range("A33").Value = "Japon"
$ws.Range("B33").Value = 14088
$ws.Range("C33").Value = 0
$ws.Range("D33").Value = 2460
$ws.Range("E33").Value = 11198
$ws.Range("F33").Value = 308
$ws.Range("G33").Value = 0
$ws.Range("H33").Value = 430

# Row 55 used to be Argentina, now shows Marruecos with its refreshed numbers
$ws.Range("A55").Value = "Marruecos"
$ws.Range("B55").Value = 4529
$ws.Range("C55").Value = 106
$ws.Range("D55").Value = 1055
$ws.Range("E55").Value = 3303
$ws.Range("F55").Value = 1
$ws.Range("G55").Value = 1
$ws.Range("H55").Value = 171

# Row 56 used to be Marruecos, now shows Argentina (its previously unchanged numbers)
$ws.Range("A56").Value = "Argentina"
$ws.Range("B56").Value = 4428
$ws.Range("C56").Value = 0
$ws.Range("D56").Value = 1256
$ws.Range("E56").Value = 2954
$ws.Range("F56").Value = 157
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 218

# --- Rows updated in place (same country, refreshed figures) ---

# Row 37: Rumania
$ws.Range("B37").Value = 12567
$ws.Range("C37").Value = 327
$ws.Range("D37").Value = 4328
$ws.Range("E37").Value = 7513
$ws.Range("F37").Value = 249

# Row 41: Dinamarca
$ws.Range("B41").Value = 9311
$ws.Range("C41").Value = 153
$ws.Range("E41").Value = 2313

# Row 67: Afganistan
$ws.Range("B67").Value = 2335
$ws.Range("C67").Value = 164
$ws.Range("D67").Value = 310
$ws.Range("E67").Value = 1957
$ws.Range("G67").Value = 4
$ws.Range("H67").Value = 68

# Row 72: Uzbekistan
$ws.Range("D72").Value = 1159
$ws.Range("E72").Value = 878

# Row 99: Libano
$ws.Range("B99").Value = 729
$ws.Range("C99").Value = 4
$ws.Range("D99").Value = 192
$ws.Range("E99").Value = 513

# Row 170: Macao
$ws.Range("D170").Value = 37
$ws.Range("E170").Value = 8
